$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")

$ws.Range("B3").Value = "Hand Seal"
$ws.Range("C3").Value = "Projectile"
$ws.Range("D3").Value = "Binding"
$ws.Range("E3").Value = "None"
$ws.Range("F3").Value = "None"

$ws.Range("I4").Select()
